# Auto-generated Excel COM-interop script to update cryptos list
# Commit: Updated cryptos list on Thu Mar 14 03:59:39 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.221.55"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.65%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.993.35"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.19%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "615.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +14.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.21"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +11.98%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.685"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.96%  "

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.85%  "

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.42%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.49"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.68%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000333"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.52%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.630.76"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.996.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.30%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.62%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.21"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.92%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.60"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.120.88"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.48%  "

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.47"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.70%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +15.95%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.99"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.10%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.28%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.55%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.09"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.86%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.21"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.52%  "

# Row 28
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.94"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.51"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.82%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.10"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.79%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.83"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.05%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.76"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.42%  "

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.56%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.54"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.91%  "

# Row 35
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "72.28"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +8.72%  "

# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +15.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "638.07"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.85%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.437"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.86%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.43"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.66%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.147"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.13%  "

# Row 41
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.12"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.82%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.07%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.57%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +4.09%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.71%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.89"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +31.56%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.862.37"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.03"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.88%  "
